$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = '@'
    $r.Value = $val
    $r.Style = 'Normal'
}

Set-TextValue 'D2' '34.779.10'
Set-TextValue 'E2' '  -0.98%  '
Set-TextValue 'D3' '1.827.46'
Set-TextValue 'E3' '  +0.47%  '
Set-TextValue 'E4' '  +0.25%  '
Set-TextValue 'D5' '230.26'
Set-TextValue 'E5' '  -1.46%  '
Set-TextValue 'E6' '  +0.37%  '
Set-TextValue 'E7' '  +0.22%  '
Set-TextValue 'D8' '39.58'
Set-TextValue 'E8' '  -3.19%  '
Set-TextValue 'D9' '0.325'
Set-TextValue 'E9' '  -0.75%  '
Set-TextValue 'D10' '0.0682'
Set-TextValue 'E10' '  -0.76%  '
Set-TextValue 'E11' '  -1.38%  '
Set-TextValue 'E12' '  +0.58%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.843.99'
Set-TextValue 'E13' '  +1.33%  '
Set-TextValue 'B14' 'Chainlink'
Set-TextValue 'C14' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D14' '11.30'
Set-TextValue 'E14' '  +1.60%  '
Set-TextValue 'D15' '0.668'
Set-TextValue 'E15' '  +0.89%  '
Set-TextValue 'E16' '  -1.41%  '
Set-TextValue 'D17' '34.780.58'
Set-TextValue 'E17' '  -0.81%  '
Set-TextValue 'D18' '69.45'
Set-TextValue 'E18' '  -0.31%  '
Set-TextValue 'E19' '  -1.03%  '
Set-TextValue 'D20' '240.06'
Set-TextValue 'E20' '  +0.04%  '
Set-TextValue 'D21' '12.11'
Set-TextValue 'E21' '  +1.94%  '
Set-TextValue 'D22' '4.66'
Set-TextValue 'E22' '  -0.75%  '
Set-TextValue 'E23' '  +0.29%  '
Set-TextValue 'D24' '2.25'
Set-TextValue 'E24' '  -0.32%  '
Set-TextValue 'D25' '171.70'
Set-TextValue 'E25' '  -0.71%  '
Set-TextValue 'D26' '7.74'
Set-TextValue 'E26' '  -1.89%  '
Set-TextValue 'E27' '  +1.95%  '
Set-TextValue 'D28' '17.30'
Set-TextValue 'E28' '  -1.39%  '
Set-TextValue 'D29' '1.50'
Set-TextValue 'E29' '  -7.87%  '
Set-TextValue 'E30' '  +0.27%  '
Set-TextValue 'D31' '0.0549'
Set-TextValue 'E31' '  -1.25%  '
Set-TextValue 'E32' '  -3.74%  '
Set-TextValue 'E33' '  -1.86%  '
Set-TextValue 'D34' '1.84'
Set-TextValue 'E34' '  +3.30%  '
Set-TextValue 'E35' '  +6.75%  '
Set-TextValue 'E36' '  +11.62%  '
Set-TextValue 'D37' '0.696'
Set-TextValue 'E37' '  +1.68%  '
Set-TextValue 'D38' '91.01'
Set-TextValue 'E38' '  -2.55%  '
Set-TextValue 'E39' '  +5.88%  '
Set-TextValue 'D40' '1.337.04'
Set-TextValue 'E40' '  +1.69%  '
Set-TextValue 'D41' '0.0192'
Set-TextValue 'E41' '  -1.20%  '
Set-TextValue 'E42' '  -1.78%  '
Set-TextValue 'D43' '2.41'
Set-TextValue 'E43' '  -2.13%  '
Set-TextValue 'E44' '  -3.36%  '
Set-TextValue 'E45' '  -0.75%  '
Set-TextValue 'D46' '6.24'
Set-TextValue 'E46' '  -2.01%  '
Set-TextValue 'E47' '  +2.30%  '
Set-TextValue 'E48' '  +0.64%  '
Set-TextValue 'E49' '  +0.19%  '
Set-TextValue 'E50' '  +3.47%  '
Set-TextValue 'E51' '  +12.90%  '

Write-Host "Applied 81 cell updates"
